$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.626.05"
$ws.Range("E2").Value = "  -0.31%  "
$ws.Range("D3").Value = "'2.087.41"
$ws.Range("E3").Value = "  +0.46%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'234.00"
$ws.Range("E5").Value = "  +0.08%  "
$ws.Range("E6").Value = "  +2.23%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "'58.14"
$ws.Range("E8").Value = "  -0.18%  "
$ws.Range("D9").Value = "'0.393"
$ws.Range("E9").Value = "  -0.81%  "
$ws.Range("D10").Value = "'0.0782"
$ws.Range("E10").Value = "  -0.25%  "
$ws.Range("E11").Value = "  +3.27%  "
$ws.Range("E12").Value = "  +2.47%  "
$ws.Range("D13").Value = "'2.396.30"
$ws.Range("E13").Value = "  +0.50%  "
$ws.Range("D14").Value = "'21.18"
$ws.Range("E14").Value = "  +1.59%  "
$ws.Range("D15").Value = "'0.775"
$ws.Range("E15").Value = "  +0.15%  "
$ws.Range("D16").Value = "'5.37"
$ws.Range("E16").Value = "  +1.42%  "
$ws.Range("D17").Value = "'2.092.50"
$ws.Range("E17").Value = "  +1.37%  "
$ws.Range("D18").Value = "'37.608.24"
$ws.Range("E18").Value = "  -0.23%  "
$ws.Range("E19").Value = "  -2.07%  "
$ws.Range("D20").Value = "'70.87"
$ws.Range("E20").Value = "  -0.38%  "
$ws.Range("D21").Value = "'0.0₃0836"
$ws.Range("E21").Value = "  +0.26%  "
$ws.Range("D22").Value = "'229.56"
$ws.Range("E22").Value = "  +0.42%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").Value = "'2.36"
$ws.Range("E24").Value = "  -1.51%  "
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").Value = "'2.39"
$ws.Range("E25").Value = "  -0.59%  "
$ws.Range("D26").Value = "'9.69"
$ws.Range("E26").Value = "  +7.17%  "
$ws.Range("D27").Value = "'170.86"
$ws.Range("E27").Value = "  +0.20%  "
$ws.Range("E28").Value = "  -3.62%  "
$ws.Range("D29").Value = "'19.54"
$ws.Range("E29").Value = "  +0.45%  "
$ws.Range("E30").Value = "  +0.39%  "
$ws.Range("E31").Value = "  +1.32%  "
$ws.Range("D32").Value = "'4.67"
$ws.Range("E32").Value = "  +0.26%  "
$ws.Range("E33").Value = "  +1.52%  "
$ws.Range("D34").Value = "'4.68"
$ws.Range("E34").Value = "  +1.09%  "
$ws.Range("E35").Value = "  +0.62%  "
$ws.Range("D36").Value = "'1.81"
$ws.Range("E36").Value = "  -0.87%  "
$ws.Range("D37").Value = "'3.33"
$ws.Range("E37").Value = "  -1.55%  "
$ws.Range("E38").Value = "  -0.04%  "
$ws.Range("D39").Value = "'5.39"
$ws.Range("E39").Value = "  +1.84%  "
$ws.Range("E40").Value = "  +8.85%  "
$ws.Range("D41").Value = "'100.97"
$ws.Range("E41").Value = "  +2.93%  "
$ws.Range("B42").Value = "Cronos"
$ws.Range("C42").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D42").Value = "'0.0960"
$ws.Range("E42").Value = "  -1.05%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "'1.20"
$ws.Range("E43").Value = "  +4.16%  "
$ws.Range("D44").Value = "'2.90"
$ws.Range("E44").Value = "  +0.71%  "
$ws.Range("D45").Value = "'16.92"
$ws.Range("E45").Value = "  +2.73%  "
$ws.Range("D46").Value = "'1.463.65"
$ws.Range("E46").Value = "  +1.00%  "
$ws.Range("E47").Value = "  -0.57%  "
$ws.Range("D48").Value = "'4.01"
$ws.Range("E48").Value = "  -5.61%  "
$ws.Range("D49").Value = "'7.26"
$ws.Range("E49").Value = "  -2.03%  "
$ws.Range("D50").Value = "'2.96"
$ws.Range("E50").Value = "  -1.70%  "
$ws.Range("D51").Value = "'2.279.66"
$ws.Range("E51").Value = "  +0.46%  "
